# Report progress devlog day 4 - tambah asset dan uji coba
#
# 1) "DEVLOG DAY 3" was previously split across two runs ("DEVLOG DAY " and
#    "3"); re-typing the same text over itself makes Word collapse it back
#    into a single run, matching the target markup.
# 2) A new "DEVLOG DAY 4" devlog entry (Heading 1 title, Heading 2 "Report"
#    sub-heading, and a Normal report paragraph) is inserted where the
#    document's trailing blank paragraph used to be, right after the
#    "DEVLOG DAY 3" entry's report text.

$d = $word.ActiveDocument

# --- Step 1: merge the two runs "DEVLOG DAY " + "3" into a single run ---
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("DEVLOG DAY 3", $true, $false, $false, $false, $false, $true, 1, $false, "DEVLOG DAY 3", 2) | Out-Null

# --- Step 2: find the (empty) paragraph right after the "Day 3" report text ---
$day3Report = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*menyambungkan pulau ke pulau*") {
        $day3Report = $cand
        break
    }
}
if ($day3Report -eq $null) {
    throw "Could not locate the 'DEVLOG DAY 3' report paragraph"
}
$target = $day3Report.Next()

# --- Fill that paragraph in with the "DEVLOG DAY 4" heading, the "Report"
#     sub-heading, and the new report paragraph, by injecting the exact
#     WordprocessingML for all three paragraphs in one shot. ---
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading1"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="0"/>
              </w:numPr>
              <w:ind w:left="360"/>
              <w:rPr>
                <w:lang w:val="en-ID"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-ID"/>
              </w:rPr>
              <w:t>DEVLOG DAY 4</w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-ID"/>
              </w:rPr>
              <w:br/>
              <w:t>SALT STUDIO CODELABS</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Heading2"/>
              <w:rPr>
                <w:lang w:val="en-ID" w:eastAsia="en-ID"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-ID" w:eastAsia="en-ID"/>
              </w:rPr>
              <w:t>Report</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:rPr>
                <w:lang w:val="en-ID" w:eastAsia="en-ID"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-ID" w:eastAsia="en-ID"/>
              </w:rPr>
              <w:t xml:space="preserve">Hari ini tim kami </w:t>
            </w:r>
            <w:r>
              <w:rPr>
                <w:lang w:val="en-ID" w:eastAsia="en-ID"/>
              </w:rPr>
              <w:t>telah selesai membuat asset untuk kucing. Kami juga telah menguploadnya ke github dan mulai mencobanya untuk digerakan. Kami juga memasukan scene falling dimana scene ini akan ditrigger apabila anak kucing terjatuh ke air. Kami memastikan dan menguji coba asset yang kami masukan. Kedepannya kami mulai melengkapi beberapa bagian yang memerlukan asset dan menguji coba terus project yang kami buat.</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$target.Range.InsertXML($xml) | Out-Null
